$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename rain model constants (shared strings used as labels in column A)
$ws.Range("A19").Value = "fi_lidar_rain_reflectivity"
$ws.Range("A20").Value = "fi_lidar_rain_intensity"

# Update the active selection to A20 (matches the saved selection state in the diff)
$ws.Range("A20").Select()
